# Natmi following Dr Hou advice
# Rebuild the LR-pairs table: rows now cover both "FAPs" and "sCs" as sending
# clusters (6 rows before -> 10 rows after), with recalculated statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# followed by the 16 numeric values for columns E..T
$data = @(
    @(2,  "FAPs", "Rspo2", "Znrf3", "ECs",  3, 1, 1.526851, 4.580553, 0.9910541119453039, 0.9910541119453039, 3, 1, 1.036231666666667, 3.108695, 0.09018629024029708, 0.09018629024029708, 1.582171356481667, 14.239542208335, 0.08937949378373905, 0.08937949378373905),
    @(3,  "FAPs", "Rspo2", "Znrf3", "FAPs", 3, 1, 1.526851, 4.580553, 0.9910541119453039, 0.9910541119453039, 3, 1, 1.678092, 5.034276, 0.1460492832155492, 0.1460492832155492, 2.562196448292, 23.059768034628, 0.1447427426774343, 0.1447427426774343),
    @(4,  "FAPs", "Rspo2", "Znrf3", "M1",   3, 1, 1.526851, 4.580553, 0.9910541119453039, 0.9910541119453039, 3, 1, 2.826506333333333, 8.479519, 0.24599916094442, 0.24599916094442, 4.315654021556333, 38.840886194007, 0.243798479989062, 0.243798479989062),
    @(5,  "FAPs", "Rspo2", "Znrf3", "M2",   3, 1, 1.526851, 4.580553, 0.9910541119453039, 0.9910541119453039, 3, 1, 2.890749333333333, 8.672248, 0.2515904182185245, 0.2515904182185245, 4.413743510349334, 39.723691593144, 0.2493397185015074, 0.2493397185015074),
    @(6,  "FAPs", "Rspo2", "Znrf3", "sCs",  3, 1, 1.526851, 4.580553, 0.9910541119453039, 0.9910541119453039, 3, 1, 3.058323, 9.174969000000001, 0.2661748473812093, 0.2661748473812093, 4.669603530873, 42.02643177785701, 0.2637936769935612, 0.2637936769935612),
    @(7,  "sCs",  "Rspo2", "Znrf3", "ECs",  1, 0.3333333333333333, 0.01378233333333333, 0.041347, 0.008945888054696122, 0.00894588805469612, 3, 1, 1.036231666666667, 3.108695, 0.09018629024029708, 0.09018629024029708, 0.01428169024055556, 0.128535212165, 0.000806796456558031, 0.0008067964565580309),
    @(8,  "sCs",  "Rspo2", "Znrf3", "FAPs", 1, 0.3333333333333333, 0.01378233333333333, 0.041347, 0.008945888054696122, 0.00894588805469612, 3, 1, 1.678092, 5.034276, 0.1460492832155492, 0.1460492832155492, 0.023128023308, 0.208152209772, 0.001306540538114913, 0.001306540538114913),
    @(9,  "sCs",  "Rspo2", "Znrf3", "M1",   1, 0.3333333333333333, 0.01378233333333333, 0.041347, 0.008945888054696122, 0.00894588805469612, 3, 1, 2.826506333333333, 8.479519, 0.24599916094442, 0.24599916094442, 0.03895585245477778, 0.350602672093, 0.002200680955357955, 0.002200680955357955),
    @(10, "sCs",  "Rspo2", "Znrf3", "M2",   1, 0.3333333333333333, 0.01378233333333333, 0.041347, 0.008945888054696122, 0.00894588805469612, 3, 1, 2.890749333333333, 8.672248, 0.2515904182185245, 0.2515904182185245, 0.03984127089511112, 0.358571438056, 0.002250699717017099, 0.002250699717017099),
    @(11, "sCs",  "Rspo2", "Znrf3", "sCs",  1, 0.3333333333333333, 0.01378233333333333, 0.041347, 0.008945888054696122, 0.00894588805469612, 3, 1, 3.058323, 9.174969000000001, 0.2661748473812093, 0.2661748473812093, 0.042150827027, 0.379357443243, 0.002381170387648123, 0.002381170387648123)
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $data) {
    $r = $row[0]
    for ($i = 1; $i -lt $row.Length; $i++) {
        $col = $columns[$i - 1]
        $ws.Range("$col$r").Value = $row[$i]
    }
}
